$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 2 (H) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 189
$wsOff.Range("C2").Value = 143
$wsOff.Range("D2").Value = 49
$wsOff.Range("E2").Value = 27

# Sheet "DEF" - row 2 (H) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 211
$wsDef.Range("C2").Value = 149
$wsDef.Range("D2").Value = 50
$wsDef.Range("E2").Value = 23
